$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 107: H107,J107,L107,N107
$ws.Range("H107").Value = 27129.79
$ws.Range("J107").Value = 807.8
$ws.Range("L107").Value = 807.8
$ws.Range("N107").Value = -4647.8

# Row 138: H138,J138,L138,N138
$ws.Range("H138").Value = 4125.294
$ws.Range("J138").Value = 4141.4526
$ws.Range("L138").Value = 12424.3578
$ws.Range("N138").Value = -22704.3578

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 28: H28,I28,K28,M28
$ws.Range("H28").Value = 5100
$ws.Range("I28").Value = 5100
$ws.Range("K28").Value = 5100
$ws.Range("M28").Value = -4908

# Row 32: H32,I32,J32,K32,L32,M32,N32
$ws.Range("H32").Value = 8154.914
$ws.Range("I32").Value = 2056.446
$ws.Range("J32").Value = 22312.072
$ws.Range("K32").Value = 2056.446
$ws.Range("L32").Value = 22312.072
$ws.Range("M32").Value = -1769.446
$ws.Range("N32").Value = -22886.072

# Row 45: H45,I45,J45,K45,L45,M45,N45
$ws.Range("H45").Value = 2524.182
$ws.Range("I45").Value = 2185.9473
$ws.Range("J45").Value = 4666.3335
$ws.Range("K45").Value = 2185.9473
$ws.Range("L45").Value = 4666.3335
$ws.Range("M45").Value = -1808.9473
$ws.Range("N45").Value = -5420.3335

# Row 62: H62,I62,J62,K62,L62,N62,M62
$ws.Range("H62").Value = 40226
$ws.Range("I62").Value = 40226
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 40226
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = -39602

# Row 65: H65,I65,J65,K65,L65,N65,M65
$ws.Range("H65").Value = 40226
$ws.Range("I65").Value = 40226
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 120678
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("M65").Value = -117558

# Row 99: H99,I99,K99,M99
$ws.Range("H99").Value = 5100
$ws.Range("I99").Value = 5100
$ws.Range("K99").Value = 5100
$ws.Range("M99").Value = -2105

# Row 110: H110,I110,K110,M110
$ws.Range("H110").Value = 12772
$ws.Range("I110").Value = 12772
$ws.Range("K110").Value = 12772
$ws.Range("M110").Value = -10727

# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 827
$ws.Range("I122").Value = 637.8
$ws.Range("K122").Value = 1913.4
$ws.Range("M122").Value = 536.6000000000001

# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 5760544
$ws.Range("I132").Value = 1511.7693
$ws.Range("J132").Value = 13247286
$ws.Range("K132").Value = 4535.3079
$ws.Range("L132").Value = 39741858
$ws.Range("M132").Value = -2005.3079
$ws.Range("N132").Value = -39746918

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 128: H128,I128,K128,M128
$ws.Range("H128").Value = 7974.25
$ws.Range("I128").Value = 7974.25
$ws.Range("K128").Value = 23922.75
$ws.Range("M128").Value = -21432.75

# Row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 9166.621999999999
$ws.Range("I134").Value = 4149.237
$ws.Range("J134").Value = 36403.855
$ws.Range("K134").Value = 12447.711
$ws.Range("L134").Value = 109211.565
$ws.Range("M134").Value = -9912.710999999999
$ws.Range("N134").Value = -114281.565

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 46535.395
$ws.Range("I31").Value = 78947.08
$ws.Range("J31").Value = 18445.268
$ws.Range("K31").Value = 78947.08
$ws.Range("L31").Value = 18445.268
$ws.Range("M31").Value = -78652.08
$ws.Range("N31").Value = -19035.268

# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 46535.395
$ws.Range("I34").Value = 78947.08
$ws.Range("J34").Value = 18445.268
$ws.Range("K34").Value = 78947.08
$ws.Range("L34").Value = 18445.268
$ws.Range("M34").Value = -78745.08
$ws.Range("N34").Value = -18849.268

# Row 74: H74,J74,L74,N74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77: H77,J77,L77,N77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 68: H68,I68,J68,K68,L68,M68,N68
$ws.Range("H68").Value = 2353.1177
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 2437.6875
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 7313.0625
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -8935.0625

# Row 71: H71,I71,J71,K71,L71,M71,N71
$ws.Range("H71").Value = 2353.1177
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 2437.6875
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 21939.1875
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -30051.1875

# Row 131: H131,I131,J131,K131,L131,N131,M131
$ws.Range("H131").Value = 1476.4
$ws.Range("I131").Value = 30
$ws.Range("J131").Value = 1491.0101
$ws.Range("K131").Value = 90
$ws.Range("L131").Value = 4473.0303
$ws.Range("N131").Value = -14553.0303
$ws.Range("M131").Value = 4950

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 102: H102,I102,J102,K102,L102,M102,N102
$ws.Range("H102").Value = 3483.0476
$ws.Range("I102").Value = 3570.7368
$ws.Range("J102").Value = 2650
$ws.Range("K102").Value = 3570.7368
$ws.Range("L102").Value = 2650
$ws.Range("M102").Value = -1948.7368
$ws.Range("N102").Value = -5894

# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 742.1
$ws.Range("I107").Value = 1047.6154
$ws.Range("K107").Value = 1047.6154
$ws.Range("M107").Value = 872.3846000000001

# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 3713.5293
$ws.Range("I126").Value = 3578.923
$ws.Range("J126").Value = 4151
$ws.Range("K126").Value = 10736.769
$ws.Range("L126").Value = 12453
$ws.Range("M126").Value = -8266.769
$ws.Range("N126").Value = -17393

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7: H7,I7,K7,M7
$ws.Range("H7").Value = 13434.875
$ws.Range("I7").Value = 18869.75
$ws.Range("K7").Value = 18869.75
$ws.Range("M7").Value = -18757.75

# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 5489.6387
$ws.Range("I122").Value = 4369.8
$ws.Range("J122").Value = 6889.4375
$ws.Range("K122").Value = 13109.4
$ws.Range("L122").Value = 20668.3125
$ws.Range("M122").Value = -10659.4
$ws.Range("N122").Value = -25568.3125

# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 13434.875
$ws.Range("I126").Value = 18869.75
$ws.Range("K126").Value = 56609.25
$ws.Range("M126").Value = -54139.25

# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1402960
$ws.Range("I132").Value = 4127.3213
$ws.Range("J132").Value = 3183292.5
$ws.Range("K132").Value = 12381.9639
$ws.Range("L132").Value = 9549877.5
$ws.Range("M132").Value = -9851.963899999999
$ws.Range("N132").Value = -9554937.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62: H62,I62,J62,K62,L62,M62,N62
$ws.Range("H62").Value = 27248.166
$ws.Range("I62").Value = 5996.3335
$ws.Range("J62").Value = 34332.11
$ws.Range("K62").Value = 5996.3335
$ws.Range("L62").Value = 34332.11
$ws.Range("M62").Value = -5372.3335
$ws.Range("N62").Value = -35580.11

# Row 65: H65,I65,J65,K65,L65,M65,N65
$ws.Range("H65").Value = 27248.166
$ws.Range("I65").Value = 5996.3335
$ws.Range("J65").Value = 34332.11
$ws.Range("K65").Value = 29981.6675
$ws.Range("L65").Value = 171660.55
$ws.Range("M65").Value = -26861.6675
$ws.Range("N65").Value = -177900.55

# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 74061.39999999999
$ws.Range("I126").Value = 116067.664
$ws.Range("J126").Value = 11052
$ws.Range("K126").Value = 348202.992
$ws.Range("L126").Value = 33156
$ws.Range("M126").Value = -345732.992
$ws.Range("N126").Value = -38096

# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1337436.4
$ws.Range("I132").Value = 6730.3613
$ws.Range("J132").Value = 7325613
$ws.Range("K132").Value = 20191.0839
$ws.Range("L132").Value = 21976839
$ws.Range("M132").Value = -17661.0839
$ws.Range("N132").Value = -21981899
